$wb = $excel.ActiveWorkbook

$productBacklog = $wb.Worksheets.Item("Product Backlog")
$productBacklog.Range("B2").Value = "Team Move Detroit To Ohio"
$productBacklog.Range("B5").Value = "Rishabh Mediratta"
$productBacklog.Range("H5").Value = "RM"
$productBacklog.Range("I5").Value = 1002035684

$productBacklog.Range("F24").Value = 1
$productBacklog.Range("G24").Value = "Finished in Sprint 1"
$productBacklog.Range("F25").Value = 1
$productBacklog.Range("G25").Value = "Finished in Sprint 1"
$productBacklog.Range("F26").Value = 1
$productBacklog.Range("G26").Value = "Finished in Sprint 1"
$productBacklog.Range("F27").Value = 1
$productBacklog.Range("G27").Value = "Finished in Sprint 1"
$productBacklog.Range("F28").Value = 2
$productBacklog.Range("G28").Value = "Finished in Sprint 2"
$productBacklog.Range("F29").Value = 2
$productBacklog.Range("G29").Value = "Finished in Sprint 2"

$sprint1 = $wb.Worksheets.Item("Sprint 01 Backlog")
$sprint1.Range("B17").Value = "CF"
$sprint1.Range("C17").Value = "RM"
$sprint1.Range("D17").Value = "Create IceCreamFlavor class extending Item and add constructor"
$sprint1.Range("E17").Value = "Completed Day 6"

$sprint1.Range("B18").Value = "MXF"
$sprint1.Range("C18").Value = "RM"
$sprint1.Range("D18").Value = "Create MixInFlavor class extending Item and add constructor"
$sprint1.Range("E18").Value = "Completed Day 6"

$sprint1.Range("B19").Value = "MX"
$sprint1.Range("C19").Value = "RM"
$sprint1.Range("D19").Value = "Create MixIn and add constructor and toString method"
$sprint1.Range("E19").Value = "Completed Day 6"

$sprint1.Range("B20").Value = "SCP"
$sprint1.Range("C20").Value = "RM"
$sprint1.Range("D20").Value = "Create Scoop class and add constructor, addMixIn and override toString methods"
$sprint1.Range("E20").Value = "Completed Day 6"

$sprint1.Range("B21").Value = "SCP"
$sprint1.Range("C21").Value = "RM"
$sprint1.Range("D21").Value = "Add TestScoop to test Scoop"
$sprint1.Range("E21").Value = "Completed Day 6"

$sprint1.Range("B22").Value = "SCP"
$sprint1.Range("C22").Value = "RM"
$sprint1.Range("D22").Value = "Fix conditions to check if mixins ArrayList is not initialized"
$sprint1.Range("E22").Value = "Completed Day 6"

$sprint2 = $wb.Worksheets.Item("Sprint 02 Backlog")
$sprint2.Range("B17").Value = "GUI"
$sprint2.Range("C17").Value = "RM"
$sprint2.Range("D17").Value = "Made a GUI instead of using a terminal"
$sprint2.Range("E17").Value = "Completed Day 6"

$sprint2.Range("B18").Value = "IGUI"
$sprint2.Range("C18").Value = "RM"
$sprint2.Range("D18").Value = "Added logic to GUI to make items and scoops"
$sprint2.Range("E18").Value = "Completed Day 6"
